$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# Numeric-looking text values in column D need to be forced to Text format
# first so Excel does not reinterpret them as numbers (e.g. "9.10" -> 9.1).

$ws.Range('D2').Value = '43.153.66'
$ws.Range('E2').Value = '  +2.55%  '
$ws.Range('D3').Value = '2.292.28'
$ws.Range('E3').Value = '  +3.48%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.34'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.04'
$ws.Range('E7').Value = '  +8.68%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.645'
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.21'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0977'
$ws.Range('E11').Value = '  +4.12%  '
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.41'
$ws.Range('E13').Value = '  +5.18%  '
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').Value = '2.635.85'
$ws.Range('E15').Value = '  +3.48%  '
$ws.Range('E16').Value = '  +5.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.871'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '2.292.47'
$ws.Range('E18').Value = '  +2.08%  '
$ws.Range('D19').Value = '43.024.53'
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('E20').Value = '  +4.59%  '
$ws.Range('E21').Value = '  +2.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.42'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.48'
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('E24').Value = '  +9.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.90'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.56'
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('E30').Value = '  -0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.96'
$ws.Range('E31').Value = '  +0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.03'
$ws.Range('E32').Value = '  +2.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.48'
$ws.Range('E33').Value = '  +6.93%  '
$ws.Range('E34').Value = '  +5.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0815'
$ws.Range('E35').Value = '  +4.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.58'
$ws.Range('E36').Value = '  +18.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.125'
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.62'
$ws.Range('E38').Value = '  +12.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.76'
$ws.Range('E39').Value = '  +3.51%  '
$ws.Range('E40').Value = '  -2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.53'
$ws.Range('E41').Value = '  +20.19%  '
$ws.Range('E42').Value = '  +5.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.95'
$ws.Range('E43').Value = '  +4.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.218'
$ws.Range('E44').Value = '  +11.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.95'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.10'
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('E48').Value = '  +3.44%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '99.06'
$ws.Range('E51').Value = '  +6.53%  '
